# Add a few missing in and outputs, try to join up pins
#
# Row 10 (Terminal block, GCT TBC05-04-1-G-G) previously only documented
# designator "U9" with a qty of 1. Two more positions (U11, U12) are now
# using this same terminal block, so the designator list and quantity are
# updated to reflect all three pin headers being tied together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Designator column (C) - join up U9 with the newly added U11, U12
$ws.Range("C10").Value = "U9, U11, U12"

# Qty column (D) - now 3 instead of 1
$ws.Range("D10").Value = 3

# Leave the selection where the edit was made, same as the saved file
$ws.Range("C11").Select()
